$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = "<10 Manufacturing" + [char]10 + "<5 Other businesses"
$ws.Range("C24").Value = "< N$ 500,000 Manufacturing" + [char]10 + "< N$ 100,000 Other businesses"
$ws.Range("D24").Value = "< N$ 1,000,000 Manufacturing" + [char]10 + "< N$ 250,000 Other businesses"

$ws.Range("B25").Value = "> 10 Manufacturing" + [char]10 + "> 5 Other businesses"
$ws.Range("C25").Value = "> N$ 500,000 Manufacturing" + [char]10 + "> N$ 100,000 Other businesses"
$ws.Range("D25").Value = "> N$ 1,000,000 Manufacturing" + [char]10 + "> N$ 250,000 Other businesses"
